$d = $word.ActiveDocument

$old0 = "Lenses are essential tools in our daily lives, from eyeglasses to telescopes, and this lesson will explore their basic principles and how they manipulate light."
$new0 = "Lenses are essential components in optical instruments like cameras, telescopes, and microscopes. They use refraction, the bending of light as it passes from one medium to another, to manipulate light and create images. "
$result0 = $d.Content.Find.Execute($old0, $true, $false, $false, $false, $false, $true, 1, $false, $new0, 2)
Write-Output "Replace 0: $result0"

$old1 = "LO-1: Students will be able to define the term 'lens' and identify different types of lenses, such as convex and concave."
$new1 = "LO-1: Students will be able to define a lens and differentiate between convex and concave lenses. "
$result1 = $d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)
Write-Output "Replace 1: $result1"

$old2 = "LO-2: Students will be able to explain how lenses refract light and demonstrate this understanding through simple experiments."
$new2 = "LO-2: Students will be able to explain how lenses refract light to form images, identifying the key concepts of focal point and focal length."
$result2 = $d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2)
Write-Output "Replace 2: $result2"

$old3 = "LO-3: Students will be able to describe real-world applications of lenses, such as in cameras, microscopes, and telescopes. "
$new3 = "LO-3: Students will be able to describe the different types of images formed by convex and concave lenses, including real and virtual, magnified and diminished, and inverted and upright. "
$result3 = $d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2)
Write-Output "Replace 3: $result3"

$old4 = "CG-1: Students will be able to define and describe the function of a lens in terms of light refraction and its application in various optical devices."
$new4 = "CG-1: Students will gain an understanding of the fundamental principles of refraction and how lenses manipulate light. "
$result4 = $d.Content.Find.Execute($old4, $true, $false, $false, $false, $false, $true, 1, $false, $new4, 2)
Write-Output "Replace 4: $result4"

$old5 = "CG-2: Students will be able to explain the relationship between the shape of a lens and its ability to converge or diverge light, understanding its impact on image formation."
$new5 = "CG-2: Students will be able to identify and explain the different types of lenses and their applications in various technologies."
$result5 = $d.Content.Find.Execute($old5, $true, $false, $false, $false, $false, $true, 1, $false, $new5, 2)
Write-Output "Replace 5: $result5"

$old6 = "CC-1: Students will be able to identify and distinguish between convex and concave lenses, recognizing their characteristics and properties."
$new6 = "CC-1: Students will be able to accurately define and explain the concept of refraction. "
$result6 = $d.Content.Find.Execute($old6, $true, $false, $false, $false, $false, $true, 1, $false, $new6, 2)
Write-Output "Replace 6: $result6"

$old7 = "CC-2: Students will be able to apply the concept of refraction to explain how lenses bend light, demonstrating their understanding through diagrams and simulations."
$new7 = "CC-2: Students will be able to analyze the behavior of light as it passes through different types of lenses."
$result7 = $d.Content.Find.Execute($old7, $true, $false, $false, $false, $false, $true, 1, $false, $new7, 2)
Write-Output "Replace 7: $result7"

$old8 = "CC-3: Students will be able to analyze the formation of images by lenses, explaining the relationship between object distance, image distance, and focal length. "
$new8 = "CC-3: Students will be able to apply their knowledge of lenses to design simple optical systems and explain their function. "
$result8 = $d.Content.Find.Execute($old8, $true, $false, $false, $false, $false, $true, 1, $false, $new8, 2)
Write-Output "Replace 8: $result8"

$old9 = "Lenses refract light.  Concave lenses diverge light. Convex lenses converge light. `n"
$new9 = "Lenses refract light to focus or spread it. Convex lenses converge light, concave lenses diverge light.  The focal length determines magnification and image size. `n"
$result9 = $d.Content.Find.Execute($old9, $true, $false, $false, $false, $false, $true, 1, $false, $new9, 2)
Write-Output "Replace 9: $result9"

$old10 = "Lenses refract light, converging or diverging rays. Convex lenses magnify, concave lenses reduce size. Lenses form images, real or virtual. `n"
$new10 = "Lenses refract light, bending it to focus or diverge.  Convex lenses converge light, forming real or virtual images. Concave lenses diverge light, creating virtual, upright images. `n"
$result10 = $d.Content.Find.Execute($old10, $true, $false, $false, $false, $false, $true, 1, $false, $new10, 2)
Write-Output "Replace 10: $result10"

$old11 = "Identify the lens type by its shape  `nFocus light using a magnifying glass  `nMeasure focal length of a convex lens `n"
$new11 = "Identify the types of lenses.  Observe light refraction through lenses. Construct a simple magnifying glass. `n"
$result11 = $d.Content.Find.Execute($old11, $true, $false, $false, $false, $false, $true, 1, $false, $new11, 2)
Write-Output "Replace 11: $result11"

$old12 = "Q-1: How do lenses manipulate light to create images, and what are the different types of lenses?"
$new12 = "Q-1: How do lenses manipulate light to create images, and what are the different types of lenses and their unique properties?"
$result12 = $d.Content.Find.Execute($old12, $true, $false, $false, $false, $false, $true, 1, $false, $new12, 2)
Write-Output "Replace 12: $result12"

$old13 = "Q-2: What are the key properties of lenses, such as focal length and magnification, and how do they affect image formation?"
$new13 = "Q-2: How does the shape and curvature of a lens influence its ability to focus light, and how can we use lenses to correct vision problems?"
$result13 = $d.Content.Find.Execute($old13, $true, $false, $false, $false, $false, $true, 1, $false, $new13, 2)
Write-Output "Replace 13: $result13"

$old14 = "Q-3: How are lenses used in everyday life and in scientific instruments, and what are the limitations of their applications? "
$new14 = "Q-3: What are the applications of lenses in various fields, such as microscopes, telescopes, cameras, and eyeglasses? "
$result14 = $d.Content.Find.Execute($old14, $true, $false, $false, $false, $false, $true, 1, $false, $new14, 2)
Write-Output "Replace 14: $result14"

$old15 = "TP-1: Lenses are curved pieces of transparent material that bend light, causing it to converge or diverge.`nTP-2:  Converging lenses bring light rays together at a focal point, forming real or virtual images depending on the object's position.`nTP-3: Diverging lenses spread light rays apart, making objects appear smaller and farther away. `n"
$new15 = "TP-1: Lenses bend light, causing it to converge or diverge, which affects how images are formed.`nTP-2: Different types of lenses, like convex and concave, have unique effects on light, leading to magnification or reduction of images.`nTP-3: Lenses are essential components in various optical instruments, including cameras, telescopes, and microscopes, enabling us to see objects that are far away or too small to be seen with the naked eye. `n"
$result15 = $d.Content.Find.Execute($old15, $true, $false, $false, $false, $false, $true, 1, $false, $new15, 2)
Write-Output "Replace 15: $result15"

$old16 = "Activity-1: Introduce lenses and their basic types: Convex and Concave. `nActivity-2: Demonstrate the effect of convex and concave lenses on light rays using a projector or a simple magnifying glass.`nActivity-3:  Encourage students to experiment with different lenses to observe how they magnify or diminish objects and create images. `n"
$new16 = "Activity-1: Explore different types of lenses (concave, convex) and their shapes. Discuss how lenses refract light.`nActivity-2: Conduct a simple experiment with a magnifying glass to demonstrate how a convex lens focuses light. `nActivity-3: Discuss real-world applications of lenses, such as eyeglasses, telescopes, and cameras. `n"
$result16 = $d.Content.Find.Execute($old16, $true, $false, $false, $false, $false, $true, 1, $false, $new16, 2)
Write-Output "Replace 16: $result16"

$old17 = "[formative assessment no 1]  What happens to the light when it passes through a convex lens? `n[formative assessment no 2]  Explain how a magnifying glass works using the concept of lenses. `n[formative assessment no 3]  If you place an object at the focal point of a converging lens, where will the image be formed? `n"
$new17 = "[formative assessment no 1] Describe how a magnifying glass uses lenses to make objects appear larger.`n[formative assessment no 2]  Explain how the shape of a lens affects the way it refracts light.`n[formative assessment no 3]  Imagine you are looking through a pair of glasses. What type of lens would be needed to correct nearsightedness and how does it work? `n"
$result17 = $d.Content.Find.Execute($old17, $true, $false, $false, $false, $false, $true, 1, $false, $new17, 2)
Write-Output "Replace 17: $result17"

$old18 = "Q-1: What is the difference between a convex and concave lens?`nQ-2: How does a magnifying glass work?`nQ-3: What happens to light when it passes through a lens? `n"
$new18 = "Q-1: What is a convex lens?`nQ-2: How does a lens bend light?`nQ-3: What is the focal length of a lens? `n"
$result18 = $d.Content.Find.Execute($old18, $true, $false, $false, $false, $false, $true, 1, $false, $new18, 2)
Write-Output "Replace 18: $result18"

$old19 = "Lenses are transparent objects that refract light, causing it to bend.  "
$new19 = "A lens is a curved piece of transparent material that refracts light, focusing or dispersing it. "
$result19 = $d.Content.Find.Execute($old19, $true, $false, $false, $false, $false, $true, 1, $false, $new19, 2)
Write-Output "Replace 19: $result19"

$old20 = "Q-1: What are the two main types of lenses?"
$new20 = "Q-1: How does a convex lens differ from a concave lens in terms of its shape and how it affects light?"
$result20 = $d.Content.Find.Execute($old20, $true, $false, $false, $false, $false, $true, 1, $false, $new20, 2)
Write-Output "Replace 20: $result20"

$old21 = "Q-2: Explain how a converging lens forms an image."
$new21 = "Q-2: Explain how a magnifying glass uses a convex lens to make objects appear larger."
$result21 = $d.Content.Find.Execute($old21, $true, $false, $false, $false, $false, $true, 1, $false, $new21, 2)
Write-Output "Replace 21: $result21"

$old22 = "Q-3: Describe how a magnifying glass works. "
$new22 = "Q-3: Design a simple experiment to demonstrate the phenomenon of refraction using a glass of water and a straw. "
$result22 = $d.Content.Find.Execute($old22, $true, $false, $false, $false, $false, $true, 1, $false, $new22, 2)
Write-Output "Replace 22: $result22"

$allResults = @($result0, $result1, $result2, $result3, $result4, $result5, $result6, $result7, $result8, $result9, $result10, $result11, $result12, $result13, $result14, $result15, $result16, $result17, $result18, $result19, $result20, $result21, $result22)
$failures = ($allResults | Where-Object { $_ -ne $true }).Count
if ($failures -gt 0) {
    throw "One or more Find/Replace operations failed ($failures failures)"
}
Write-Output "All $($allResults.Count) replacements applied successfully"
